# Scheduled runner update: refresh cached Universalis market-price derived
# columns (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across all job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly pulled values.

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 3800.1333
$ws.Range("I19").Value = 2322.2222
$ws.Range("J19").Value = 6017
$ws.Range("K19").Value = 2322.2222
$ws.Range("L19").Value = 6017
$ws.Range("M19").Value = -2147.2222
$ws.Range("N19").Value = -6367

# Row 58
$ws.Range("H58").Value = 620.2273
$ws.Range("I58").Value = 180.88889
$ws.Range("J58").Value = 924.38464
$ws.Range("K58").Value = 542.6666700000001
$ws.Range("L58").Value = 2773.15392
$ws.Range("M58").Value = -392.6666700000001
$ws.Range("N58").Value = -3073.15392

# Row 125
$ws.Range("H125").Value = 3748.6667
$ws.Range("I125").Value = 5225
$ws.Range("J125").Value = 796
$ws.Range("K125").Value = 47025
$ws.Range("L125").Value = 7164
$ws.Range("M125").Value = -44565
$ws.Range("N125").Value = -12084

# Row 127
$ws.Range("H127").Value = 469202.75
$ws.Range("I127").Value = 325
$ws.Range("J127").Value = 674336.75
$ws.Range("K127").Value = 975
$ws.Range("L127").Value = 2023010.25
$ws.Range("M127").Value = 3985
$ws.Range("N127").Value = -2032930.25

# Row 137
$ws.Range("H137").Value = 14286625
$ws.Range("I137").Value = 812.8333
$ws.Range("K137").Value = 2438.4999
$ws.Range("M137").Value = 111.5001000000002


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1852.1666
$ws.Range("I2").Value = 1475
$ws.Range("J2").Value = 2606.5
$ws.Range("K2").Value = 1475
$ws.Range("L2").Value = 2606.5
$ws.Range("M2").Value = -1362
$ws.Range("N2").Value = -2832.5

# Row 45
$ws.Range("H45").Value = 1475
$ws.Range("I45").Value = 1422.2222
$ws.Range("J45").Value = 1542.8572
$ws.Range("K45").Value = 1422.2222
$ws.Range("L45").Value = 1542.8572
$ws.Range("M45").Value = -1045.2222
$ws.Range("N45").Value = -2296.8572

# Row 63
$ws.Range("H63").Value = 20001124
$ws.Range("J63").Value = 2498
$ws.Range("L63").Value = 2498
$ws.Range("N63").Value = -3870

# Row 66
$ws.Range("H66").Value = 20001124
$ws.Range("J66").Value = 2498
$ws.Range("L66").Value = 12490
$ws.Range("N66").Value = -19354

# Row 116
$ws.Range("H116").Value = 1852.1666
$ws.Range("I116").Value = 1475
$ws.Range("J116").Value = 2606.5
$ws.Range("K116").Value = 1475
$ws.Range("L116").Value = 2606.5
$ws.Range("M116").Value = 819
$ws.Range("N116").Value = -7194.5

# Row 117
$ws.Range("H117").Value = 22983.385
$ws.Range("J117").Value = 22983.385
$ws.Range("L117").Value = 22983.385
$ws.Range("N117").Value = -32161.385

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1852.1666
$ws.Range("I3").Value = 1475
$ws.Range("J3").Value = 2606.5
$ws.Range("K3").Value = 1475
$ws.Range("L3").Value = 2606.5
$ws.Range("M3").Value = -1361
$ws.Range("N3").Value = -2834.5

# Row 75
$ws.Range("H75").Value = 40168.57
$ws.Range("I75").Value = 10000
$ws.Range("J75").Value = 45196.668
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 45196.668
$ws.Range("M75").Value = -9064
$ws.Range("N75").Value = -47068.668

# Row 78
$ws.Range("H78").Value = 40168.57
$ws.Range("I78").Value = 10000
$ws.Range("J78").Value = 45196.668
$ws.Range("K78").Value = 30000
$ws.Range("L78").Value = 135590.004
$ws.Range("M78").Value = -25320
$ws.Range("N78").Value = -144950.004


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 111
$ws.Range("H111").Value = 50702
$ws.Range("J111").Value = 50702
$ws.Range("L111").Value = 50702
$ws.Range("N111").Value = -58882


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 690.7143
$ws.Range("I5").Value = 265.45456
$ws.Range("J5").Value = 2250
$ws.Range("K5").Value = 796.36368
$ws.Range("L5").Value = 6750
$ws.Range("M5").Value = -684.36368
$ws.Range("N5").Value = -6974

# Row 40
$ws.Range("H40").Value = 331.25
$ws.Range("I40").Value = 167.5
$ws.Range("J40").Value = 495
$ws.Range("K40").Value = 670
$ws.Range("L40").Value = 1980
$ws.Range("M40").Value = -601
$ws.Range("N40").Value = -2118

# Row 87
$ws.Range("H87").Value = 15000
$ws.Range("I87").Value = 11666.667
$ws.Range("J87").Value = 20000
$ws.Range("K87").Value = 35000.001
$ws.Range("L87").Value = 60000
$ws.Range("M87").Value = -33752.001
$ws.Range("N87").Value = -62496

# Row 90
$ws.Range("H90").Value = 15000
$ws.Range("I90").Value = 11666.667
$ws.Range("J90").Value = 20000
$ws.Range("K90").Value = 105000.003
$ws.Range("L90").Value = 180000
$ws.Range("M90").Value = -98760.003
$ws.Range("N90").Value = -192480

# Row 92
$ws.Range("H92").Value = 639.6667
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 639.6667
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1919.0001
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4415.0001

# Row 113
$ws.Range("H113").Value = 805.875
$ws.Range("I113").Value = 585.8
$ws.Range("J113").Value = 963.0714
$ws.Range("K113").Value = 1757.4
$ws.Range("L113").Value = 2889.2142
$ws.Range("M113").Value = 412.6000000000001
$ws.Range("N113").Value = -7229.2142

# Row 116
$ws.Range("H116").Value = 567.1667
$ws.Range("I116").Value = 480.6
$ws.Range("K116").Value = 1441.8
$ws.Range("M116").Value = 2000.2

# Row 133
$ws.Range("H133").Value = 76926230
$ws.Range("I133").Value = 90910960
$ws.Range("J133").Value = 10250
$ws.Range("K133").Value = 272732880
$ws.Range("L133").Value = 30750
$ws.Range("M133").Value = -272727820
$ws.Range("N133").Value = -40870

# Row 135
$ws.Range("H135").Value = 690.7143
$ws.Range("I135").Value = 265.45456
$ws.Range("J135").Value = 2250
$ws.Range("K135").Value = 2389.09104
$ws.Range("L135").Value = 20250
$ws.Range("M135").Value = 145.9089599999998
$ws.Range("N135").Value = -25320

# Row 137
$ws.Range("H137").Value = 4806.8667
$ws.Range("I137").Value = 1918.8889
$ws.Range("J137").Value = 9138.833000000001
$ws.Range("K137").Value = 5756.6667
$ws.Range("L137").Value = 27416.499
$ws.Range("M137").Value = -656.6666999999998
$ws.Range("N137").Value = -37616.499

# Row 139
$ws.Range("H139").Value = 2810.5881
$ws.Range("I139").Value = 1484.2858
$ws.Range("J139").Value = 9000
$ws.Range("K139").Value = 4452.857400000001
$ws.Range("L139").Value = 27000
$ws.Range("M139").Value = 687.1425999999992
$ws.Range("N139").Value = -37280


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

# Row 107
$ws.Range("H107").Value = 1990.1428
$ws.Range("I107").Value = 3007.75
$ws.Range("J107").Value = 633.3333
$ws.Range("K107").Value = 3007.75
$ws.Range("L107").Value = 633.3333
$ws.Range("M107").Value = -1087.75
$ws.Range("N107").Value = -4473.3333

# Row 118
$ws.Range("H118").Value = 14378.261
$ws.Range("J118").Value = 14378.261
$ws.Range("L118").Value = 14378.261
$ws.Range("N118").Value = -17692.261


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5243.4062
$ws.Range("I7").Value = 5500.2666
$ws.Range("J7").Value = 5016.7646
$ws.Range("K7").Value = 5500.2666
$ws.Range("L7").Value = 5016.7646
$ws.Range("M7").Value = -5388.2666
$ws.Range("N7").Value = -5240.7646

# Row 21
$ws.Range("H21").Value = 85007
$ws.Range("J21").Value = 85007
$ws.Range("L21").Value = 85007
$ws.Range("N21").Value = -85355

# Row 22
$ws.Range("H22").Value = 1265.25
$ws.Range("I22").Value = 556
$ws.Range("J22").Value = 1419.4348
$ws.Range("K22").Value = 556
$ws.Range("L22").Value = 1419.4348
$ws.Range("M22").Value = -261
$ws.Range("N22").Value = -2009.4348

# Row 27
$ws.Range("H27").Value = 1265.25
$ws.Range("I27").Value = 556
$ws.Range("J27").Value = 1419.4348
$ws.Range("K27").Value = 556
$ws.Range("L27").Value = 1419.4348
$ws.Range("M27").Value = -449
$ws.Range("N27").Value = -1633.4348

# Row 46
$ws.Range("H46").Value = 909.9091
$ws.Range("I46").Value = 695
$ws.Range("K46").Value = 695
$ws.Range("M46").Value = -507

# Row 81
$ws.Range("H81").Value = 29900
$ws.Range("J81").Value = 39800
$ws.Range("L81").Value = 39800
$ws.Range("N81").Value = -41796

# Row 84
$ws.Range("H84").Value = 29900
$ws.Range("J84").Value = 39800
$ws.Range("L84").Value = 119400
$ws.Range("N84").Value = -129384

# Row 93
$ws.Range("H93").Value = 950.5
$ws.Range("I93").Value = 650.75
$ws.Range("K93").Value = 650.75
$ws.Range("M93").Value = 597.25

# Row 126
$ws.Range("H126").Value = 5243.4062
$ws.Range("I126").Value = 5500.2666
$ws.Range("J126").Value = 5016.7646
$ws.Range("K126").Value = 16500.7998
$ws.Range("L126").Value = 15050.2938
$ws.Range("M126").Value = -14030.7998
$ws.Range("N126").Value = -19990.2938

# Row 136
$ws.Range("H136").Value = 22728404
$ws.Range("I136").Value = 22728404
$ws.Range("K136").Value = 68185212
$ws.Range("M136").Value = -68182662


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 777
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 777
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 777
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1003

# Row 81
$ws.Range("H81").Value = 1666.6666
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = -8122

# Row 84
$ws.Range("H84").Value = 1666.6666
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = -40608

